$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 1049
$ws1.Range("F6").Value = 370
$ws1.Range("F8").Value = 569
$ws1.Range("F12").Value = 3036
$ws1.Range("F13").Value = 531
$ws1.Range("F15").Value = 1481
$ws1.Range("F16").Value = 824
$ws1.Range("F20").Value = 67
$ws1.Range("F21").Value = 1162
$ws1.Range("F22").Value = 233
$ws1.Range("F25").Value = 3608
$ws1.Range("F28").Value = 1592
$ws1.Range("F29").Value = 56

# Sheet "演出" (sheet2)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F9").Value = 33

# Sheet "全部类型" (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F14").Value = 33
$ws4.Range("F16").Value = 1049
$ws4.Range("F17").Value = 370
$ws4.Range("F19").Value = 569
$ws4.Range("F23").Value = 3036
$ws4.Range("F24").Value = 531
$ws4.Range("F26").Value = 1481
$ws4.Range("F27").Value = 824
$ws4.Range("F31").Value = 67
$ws4.Range("F34").Value = 1162
$ws4.Range("F35").Value = 233
$ws4.Range("F38").Value = 3608
$ws4.Range("F41").Value = 1592
$ws4.Range("F44").Value = 56
